$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '42.432.05'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.40%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.246.51'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.01%  '

$ws.Range("E4").Value = '  -0.02%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '246.23'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.52%  '

$ws.Range("E6").Value = '  +0.04%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '75.67'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.52%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.618'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -1.79%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '43.90'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +9.09%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0949'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.21%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '7.26'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.28%  '

$ws.Range("E13").Value = '  -1.43%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '2.585.97'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.08%  '

$ws.Range("E15").Value = '  -2.24%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.853'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -1.06%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.248.85'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.37%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '42.269.13'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.13%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.0000102'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +4.29%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '6.17'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.03%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '72.19'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.87%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '2.23'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.98%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '231.74'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '9.28'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +30.84%  '

$ws.Range("E25").Value = '  +0.08%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '11.48'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +3.24%  '

$ws.Range("E27").Value = '  -3.14%  '

$ws.Range("E28").Value = '  -0.70%  '

$ws.Range("E29").Value = '  +1.49%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '167.94'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.45%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '20.69'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.63%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.0826'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.67%  '

$ws.Range("E33").Value = '  -0.44%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '30.58'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -6.60%  '

$ws.Range("E35").Value = '  +11.20%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.126'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.32%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '4.53'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.42%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.0315'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +5.86%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '13.97'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +5.94%  '

$ws.Range("E40").Value = '  -1.47%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '5.80'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -2.87%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '64.08'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +6.42%  '

$ws.Range("E43").Value = '  -0.45%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '107.89'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -8.47%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.81'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.09%  '

$ws.Range("E46").Value = '  +1.47%  '

$ws.Range("E47").Value = '  +0.00%  '

$ws.Range("E48").Value = '  -0.07%  '

$ws.Range("E49").Value = '  +0.76%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.35'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +5.57%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '4.15'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -2.48%  '
